$d = $word.ActiveDocument

# --- Step 1: in-place text replacements via Find/Replace for paragraphs that already
#     contain text (keeps them from gaining a spurious xml:space="preserve"). ---

# Heading
$d.Content.Find.Execute("Review 158: [Short] Linguistic Binding in Diffusion Models: Enhancing Attribute Correspondence through Attention Map Alignment", $true, $false, $false, $false, $false, $true, 1, $false, "Review 157: [Short] End-to-End Speech Recognition Contextualization with Large Language Models, 30.09.23", 2) | Out-Null

# Bold "Paper:" line
$d.Content.Find.Execute("Paper: https://arxiv.org/abs/2306.08877v3", $true, $false, $false, $false, $false, $true, 1, $false, "Paper: https://arxiv.org/abs/2309.10917v1", 2) | Out-Null

# Plain link line
$d.Content.Find.Execute("https://arxiv.org/abs/2306.08877", $true, $false, $false, $false, $false, $true, 1, $false, "https://huggingface.co/papers/2309.10917", 2) | Out-Null

# "מודלי דיפוזיה..." paragraph -> "איך הם עשו זאת?..."
$d.Content.Find.Execute("מודלי דיפוזיה מודרניים מצטיינים ביצירת תמונות באיכות מרהיבה מתיאור טקסטואלי (ובטח DALLe3) וברוב המקרים התמונה ממש מתאימה לתיאור. אולם עדיין יש מקרים שמודל מתבלבל למשל בין הצבעים של האובייקטים המופיעים בתיאור. היום ב-#shorthebrewpapereviews סוקרים מאמר 🇮🇱 המציע שיטה למניעת בלבול סמנטי בין תכונות האובייקטים בתמונה. ", $true, $false, $false, $false, $false, $true, 1, $false, "איך הם עשו זאת? מכיוון שאי אפשר סתם לקחת אות אודיו להזין אותו כמו שהוא למודל שפה נדרש כאן אנקודר שמקודד את הפיצ'רים המהותיים של אות אודיו. מחברי המאמר משתמשים במודל מאומן מראש הנקרא ConFormer ומפיק לנו ייצוג לטנטי של אות אודיו (כלומר מערך של וקטורים המייצגים כל מקטע של אודיו או בפשטות טוקני אודיו). ד״א ConFormer הוא מודל די מעניין (הוצע על ידי גוגל) המשלב ארכיטקטורת הטרנספורמר עם שכבות קונבולוציה (משתמשים שם גם בקידוד מיקום יחסי RoPE שנהיה מאוד פופולרי היום).", 2) | Out-Null

# "הגישה המוצעת..." paragraph -> "לאחר מכן לוקחים..."
$d.Content.Find.Execute("הגישה המוצעת הינה פשוטה ואלגנטית. בשלב הראשון המחברים בונים את גרף התלויות הסינטקטית של הפרופמט כלומר מפיקים את כל קבוצות המילים (נגיד שם עצם ושם תואר) המתאימים אחד לשני (כמו (ארנב, צהוב) או (כורסא, בסגנון, מלון). לאחר מכן המחברים מכיילים מודל שפה עם פונקציה לוס ש״מפקחת״ על הדיוק הסמנטי של האובייקטים בתמונה. ", $true, $false, $false, $false, $false, $true, 1, $false, "לאחר מכן לוקחים את ייצוג של טוקני האודיו ומזינים אותם למודל שפה מאומן (הם לקחו LLAMA) יחד עם עוד מידע על האודיו כמו שם הוידאו שממנו הוא נלקח או התיאור הטקסטואלי. בסוף מטייבים(fine-tune) מודל שפה בסגנון LoRA על דאטהסט המורכב מזוגות של אודיו והטקסט. ולהפתעתי זה עובד ממש לא רע.", 2) | Out-Null

# --- Step 2: insert 3 empty Normal paragraphs + 1 new paragraph with the new
#     "אודיו וטקסט?" text, right after the (now-updated) link paragraph (#4). ---
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertParagraphAfter()
$d.Paragraphs.Item(8).Range.Text = "אודיו וטקסט? נכון ששני סוגי דאטה אלו הם די שונים ולא הגיוני להניח שמודל שאומן על טקסט יכול להביא תוצאות טובות גם על אודיו לאחר כיול קל. אך התברר שזה אפשרי. במאמר שנסקור היום ב-#shorthebrewpapereviews למעשה לקחו מודל שפה מאומן והשתמשו בו בשביל לבצע משימה audio2text. כלומר להפיק את מה שנאמר בקטע אודיו. "

# --- Step 3: delete the two trailing paragraphs that no longer exist in the new
#     review (the old "איך זה נעשה?..." + trailing <w:br/>, and "איבר נוסף..."). ---
#     After the 4-paragraph insertion above, these shifted from index 10/11 -> 14/15.
$pStart = $d.Paragraphs.Item(14)
$pEnd = $d.Paragraphs.Item(15)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

Write-Output ("FinalParaCount=" + $d.Paragraphs.Count)
